$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 45.8364935247954

# Row 3 values
$ws.Range("B3").Value = 0.01253208636536152
$ws.Range("C3").Value = 109.9114832445916
$ws.Range("D3").Value = 18.71679738969934
$ws.Range("E3").Value = 2459690191846.092
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2459690191974.732
